# Add the new "Save" column (H) to the header row and its data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from G1 (bold/centered/bordered header style)
# onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column.
$ws.Range("H2").Value = 0
